$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain (non-numeric-looking) text updates - safe to set directly.
$plainUpdates = @(
    @('D2', '26.848.79'),
    @('E2', '  -1.89%  '),
    @('D3', '1.802.32'),
    @('E3', '  -1.29%  '),
    @('E5', '  -1.84%  '),
    @('E6', '  -0.06%  '),
    @('E7', '  +3.90%  '),
    @('E8', '  -2.42%  '),
    @('E9', '  -1.63%  '),
    @('E10', '  -2.23%  '),
    @('E11', '  -3.34%  '),
    @('D12', '1.837.02'),
    @('E12', '  +0.51%  '),
    @('E13', '  -1.95%  '),
    @('E14', '  -1.88%  '),
    @('E15', '  -3.79%  '),
    @('E16', '  -1.30%  '),
    @('E17', '  -0.20%  '),
    @('E18', '  -1.13%  '),
    @('E19', '  -0.05%  '),
    @('E20', '  -3.36%  '),
    @('D21', '26.852.21'),
    @('E21', '  -1.93%  '),
    @('E22', '  -2.20%  '),
    @('E23', '  -3.97%  '),
    @('D24', '2.031.62'),
    @('E24', '  -1.29%  '),
    @('E25', '  -3.04%  '),
    @('E27', '  -1.84%  '),
    @('E28', '  -7.31%  '),
    @('E29', '  -3.18%  '),
    @('E30', '  -1.45%  '),
    @('E31', '  +0.06%  '),
    @('E32', '  -4.11%  '),
    @('E33', '  -4.79%  '),
    @('E34', '  +0.28%  '),
    @('E35', '  -3.41%  '),
    @('E37', '  -0.82%  '),
    @('E38', '  -2.10%  '),
    @('E39', '  -1.33%  '),
    @('E40', '  +1.82%  '),
    @('E41', '  -1.49%  '),
    @('B42', 'TheSandbox'),
    @('C42', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'),
    @('E42', '  -1.29%  '),
    @('B43', 'RenderToken'),
    @('C43', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'),
    @('E43', '  +1.70%  '),
    @('E44', '  -3.61%  '),
    @('E45', '  -2.14%  '),
    @('E46', '  -2.17%  '),
    @('E47', '  -3.29%  '),
    @('E48', '  -1.48%  '),
    @('E49', '  -0.08%  '),
    @('E50', '  -2.21%  ')
)

foreach ($u in $plainUpdates) {
    $ws.Range($u[0]).Value = $u[1]
}

# Numeric-looking text updates (e.g. "1.000", "0.9998") - Excel would
# otherwise coerce these to actual numbers, so force them in as text via a
# leading quote-prefix, then reset the cell style back to Normal so no
# stray "quote prefix" style/format gets attached to the cell.
$textNumberUpdates = @(
    @('D4', '1.000'),
    @('D6', '0.9998'),
    @('D7', '0.4651'),
    @('D8', '0.3698'),
    @('D9', '0.07357'),
    @('D10', '0.8676'),
    @('D11', '20.35'),
    @('D13', '5.355'),
    @('D14', '92.15'),
    @('D15', '6.500'),
    @('D16', '0.07026'),
    @('D17', '0.9995'),
    @('D18', '0.000008695'),
    @('D19', '0.9999'),
    @('D22', '5.289'),
    @('D23', '10.55'),
    @('D25', '1.901'),
    @('D26', '151.52'),
    @('D28', '2.138'),
    @('D29', '5.243'),
    @('D30', '116.12'),
    @('D31', '0.08897'),
    @('D32', '0.7593'),
    @('D33', '1.151'),
    @('D34', '2.933'),
    @('D35', '4.457'),
    @('D36', '0.9994'),
    @('D37', '1.101'),
    @('D38', '0.01951'),
    @('D39', '0.05243'),
    @('D40', '2.928'),
    @('D41', '7.216'),
    @('D42', '0.5286'),
    @('D43', '2.360'),
    @('D44', '0.1660'),
    @('D45', '8.496'),
    @('D46', '0.5008'),
    @('D47', '10.28'),
    @('D48', '103.95'),
    @('D49', '0.9993'),
    @('D50', '1.661'),
    @('D51', '0.06284')
)

foreach ($u in $textNumberUpdates) {
    $ws.Range($u[0]).Value = "'" + $u[1]
    $ws.Range($u[0]).Style = "Normal"
}

Write-Output "Applied $($plainUpdates.Count) plain + $($textNumberUpdates.Count) numeric-text updates"
